$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: G1 "Elapsed Time", H1 "CPU" (same bold/centered header style as A1:F1)
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Corrected B2 value (tiny float precision bump)
$ws.Range("B2").Value = 0.03244348355771106

# New data cells for row 2: elapsed time + cpu usage
$ws.Range("G2").Value = 0.1289622459000384
$ws.Range("H2").Value = 0.991
